$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# --- Data table updates ---
$ws.Range("L14").Value = 0
$ws.Range("F15").Value = "'0"
$ws.Range("H15").Value = -100
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 5.555555555555
$ws.Range("I16").Value = 104
$ws.Range("J16").Value = 126
$ws.Range("K16").Value = -17.460317460317
$ws.Range("L16").Value = 48.571428571428
$ws.Range("M16").Value = -38.095238095238
$ws.Range("N16").Value = -85.014409221902
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = 9.090909090909
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = -6.25
$ws.Range("I17").Value = 233
$ws.Range("J17").Value = 244
$ws.Range("K17").Value = -4.508196721311
$ws.Range("L17").Value = 43.827160493827
$ws.Range("M17").Value = 62.937062937062
$ws.Range("N17").Value = 8.878504672897
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 10
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 92
$ws.Range("J18").Value = 84
$ws.Range("K18").Value = 9.523809523809
$ws.Range("L18").Value = 29.577464788732
$ws.Range("M18").Value = -52.577319587628
$ws.Range("N18").Value = -89.264877479579
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = -16.666666666666
$ws.Range("I19").Value = 192
$ws.Range("J19").Value = 220
$ws.Range("K19").Value = -12.727272727272
$ws.Range("L19").Value = 38.129496402877
$ws.Range("M19").Value = -3.517587939698
$ws.Range("N19").Value = -43.026706231454
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 83.333333333333
$ws.Range("F20").Value = 31
$ws.Range("H20").Value = 24
$ws.Range("I20").Value = 173
$ws.Range("J20").Value = 158
$ws.Range("K20").Value = 9.493670886075
$ws.Range("L20").Value = 21.830985915493
$ws.Range("M20").Value = 6.134969325153
$ws.Range("N20").Value = -91.828058573453
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 3.225806451612
$ws.Range("F21").Value = 115
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = -2.542372881355
$ws.Range("I21").Value = 808
$ws.Range("J21").Value = 847
$ws.Range("K21").Value = -4.604486422668
$ws.Range("L21").Value = 33.996683250414
$ws.Range("M21").Value = -7.972665148063
$ws.Range("N21").Value = -80.95238095238
$ws.Range("C22").Value = "'0"
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = -41.666666666666
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -9.375
$ws.Range("F24").Value = 110
$ws.Range("G24").Value = 122
$ws.Range("H24").Value = -9.83606557377
$ws.Range("I24").Value = 922
$ws.Range("J24").Value = 945
$ws.Range("K24").Value = -2.433862433862
$ws.Range("L24").Value = 56.536502546689
$ws.Range("M24").Value = 92.484342379958
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = 11.111111111111
$ws.Range("I25").Value = 339
$ws.Range("J25").Value = 337
$ws.Range("K25").Value = 0.593471810089
$ws.Range("L25").Value = 18.531468531468
$ws.Range("M25").Value = -12.853470437018
$ws.Range("C26").Value = "'0"
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 100
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -14.285714285714
$ws.Range("I27").Value = 42
$ws.Range("J27").Value = 32
$ws.Range("K27").Value = 31.25
$ws.Range("L27").Value = 0
$ws.Range("L28").Value = -46.153846153846
$ws.Range("L29").Value = -44.444444444444

# --- Fix styles for cells that changed numeric<->text type ---
$ws.Range("D14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E16").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E27").PasteSpecial(-4122)
